$d = $word.ActiveDocument

# 1) Update the six existing bullet paragraphs in place using Find/Replace.
$replacements = @(
    @{
        Old = "Protocol: SIDs URNs Resources. Endpoints: Case Classes Events Signatures, Statement Data Pattern Matching Events."
        New = "Protocol: SIDs URNs Resources. Endpoints: Case Classes Aggregated Message Signatures, Aligned Statements Data Pattern Matching Message Events Resource Statement Occurrences. Resource Monad."
    },
    @{
        Old = "Core Model Upper Resources (DCI Context / Facets: Metaclass, Class, etc. as Resource, root navigation Context Resource)."
        New = "Core Model Upper Resources (DCI Context / Facets: Metaclass, Class, etc. as Resource, root navigation Context Resource). Aggregation (schema cases) / Alignment (resource statements occurrences): Activation."
    },
    @{
        Old = "Protocol: GET URN Case Classes / Statement Data Aggregated Events Messages Statements."
        New = "Core Model Functional Transforms: Functional Activation Statements:  Aggregation Schema Case Classes Statements / Alignment Message Events Resource Statement Occurrence."
    },
    @{
        Old = "Protocol: Browse Messages Events Statements. Build Context State Flows."
        New = "Protocol: GET URN Case Classes (Aggregation) / Statements Data (Alignment) Message Events Resource Statement Occurrences."
    },
    @{
        Old = "Protocol: POST URN Navigation Context built Case Class Statement Data Events."
        New = "Protocol: GET Browse Resource Aggregated / Aligned Message Events Resource Statement Occurrences. Build Context State Flows (Monad Functional Activation)."
    },
    @{
        Old = "Protocol: POST Subsequent entailed Context Browsing / Events Transforms."
        New = "Protocol: POST URN Navigation Context State Built Resource Activation Data Statements."
    }
)

foreach ($rep in $replacements) {
    $d.Content.Find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2)
}

# 2) Insert three brand-new bulleted paragraphs right after the paragraph that
#    now reads "Protocol: POST URN Navigation Context State Built Resource
#    Activation Data Statements." (the 6th updated bullet above), matching the
#    same list formatting (numId 3, ilvl 0, ind left=600 hanging=360).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Navigation Context State Built Resource Activation Data Statements*") {
        $anchor = $p
        break
    }
}

$newTexts = @(
    "Protocol: POST Subsequent entailed Context Browsing / Events Functional Transforms Activations.",
    "Monad: Resources (Metaclass, Class, etc.). Context.",
    "Transform: Statements (schema and occurrences)."
)

$prev = $anchor
foreach ($t in $newTexts) {
    $prev.Range.InsertParagraphAfter()
    $newPara = $prev.Next()
    $newPara.Range.Text = $t
    $prev = $newPara
}
